# Weekly data update for "Hortaliza, Terminal La Palmera de La Serena - Apio"
# Insert two new rows (a new "Primera" and "Segunda" entry) at the top of the
# data block (just after the existing row 196), pushing all the existing
# rows down by 2. Then fill in the two new rows with their values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at row 197 (shifts old rows 197..254 down to 199..256)
$ws.Range("A197:R198").EntireRow.Insert()

# --- New row 197 ("Primera") ---
$ws.Range("A197").Value = 8
$ws.Range("B197").Value = "Terminal La Palmera de La Serena"
$ws.Range("C197").Value = "Coquimbo"
$ws.Range("D197").Value = 44508
$ws.Range("E197").Value = 4
$ws.Range("F197").Value = 100112017
$ws.Range("G197").Value = "Apio"
$ws.Range("H197").Value = "Americana (o)"
$ws.Range("I197").Value = "Primera"
$ws.Range("J197").Value = 2400
$ws.Range("K197").Value = 7000
$ws.Range("L197").Value = 8000
$ws.Range("M197").Value = 7500
$ws.Range("N197").Value = "`$/docena de matas"
$ws.Range("O197").Value = "Provincia del Elquí"
$ws.Range("P197").Value = 1250
$ws.Range("Q197").Value = 6
$ws.Range("R197").Value = "Hortaliza"

# --- New row 198 ("Segunda") ---
$ws.Range("A198").Value = 8
$ws.Range("B198").Value = "Terminal La Palmera de La Serena"
$ws.Range("C198").Value = "Coquimbo"
$ws.Range("D198").Value = 44508
$ws.Range("E198").Value = 4
$ws.Range("F198").Value = 100112017
$ws.Range("G198").Value = "Apio"
$ws.Range("H198").Value = "Americana (o)"
$ws.Range("I198").Value = "Segunda"
$ws.Range("J198").Value = 1400
$ws.Range("K198").Value = 5500
$ws.Range("L198").Value = 6000
$ws.Range("M198").Value = 5750
$ws.Range("N198").Value = "`$/docena de matas"
$ws.Range("O198").Value = "Provincia del Elquí"
$ws.Range("P198").Value = 958
$ws.Range("Q198").Value = 6
$ws.Range("R198").Value = "Hortaliza"
